$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rebuild hyperlinks ---
# NOTE: this engine's Hyperlink object setters (Address/TextToDisplay/ScreenTip) and
# Range.Hyperlinks.Delete() operate across the *whole worksheet* rather than being
# scoped to a single cell/hyperlink, so the reliable way to update a hyperlink is to
# clear every hyperlink once and re-add all of them with their final values. The
# TextToDisplay argument also drives the cell's literal text, so cells whose text
# must end up different from the hyperlink's stored display text (A2, G2, J2, M2)
# are corrected with a plain .Value write afterwards.
$ws.Range("A1").Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("F2"), "mailto:C@bi`$ush5", [System.Type]::Missing, [System.Type]::Missing, [System.Type]::Missing)
$ws.Hyperlinks.Add($ws.Range("G2"), "https://test17.cliotest.com/cabicentral/control/main", [System.Type]::Missing, [System.Type]::Missing, "https://test17.cliotest.com/cabicentral/control/main")
$ws.Hyperlinks.Add($ws.Range("J2"), "https://test17.cliotest.com/warehouse/control/main", [System.Type]::Missing, [System.Type]::Missing, "https://test19.cliotest.com/warehouse/control/main")
$ws.Hyperlinks.Add($ws.Range("A2"), "https://test4.cliotest.com/backoffice/control/main", [System.Type]::Missing, [System.Type]::Missing, "https://test4.cliotest.com/backoffice/control/main")
$ws.Hyperlinks.Add($ws.Range("M2"), "https://mirandakate.cabitest19.com/", [System.Type]::Missing, [System.Type]::Missing, "https://mirandakate.cabitest19.com")
$ws.Hyperlinks.Add($ws.Range("N2"), "mailto:michigan@na.com", [System.Type]::Missing, [System.Type]::Missing, [System.Type]::Missing)

# --- Update the displayed cell text (test19 -> test21, cabitest19 -> cabitest21) ---
$ws.Range("A2").Value = "https://test21.cliotest.com/backoffice/control/main"
$ws.Range("G2").Value = "https://test21.cliotest.com/cabicentral/control/main"
$ws.Range("J2").Value = "https://test21.cliotest.com/warehouse/control/main"
$ws.Range("M2").Value = "https://mirandakate.cabitest21.com"

# --- Move the active selection from M2 to J2 ---
$ws.Range("J2").Select()
